# Update "想去人数" (attendance count) figures in column F across the
# three sheets that list con / event data: 展览, 演出, 全部类型.
# (本地生活 has no rows and is untouched.)

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 5662
$ws1.Range("F4").Value  = 7670
$ws1.Range("F10").Value = 40
$ws1.Range("F11").Value = 4446
$ws1.Range("F15").Value = 3001
$ws1.Range("F19").Value = 551
$ws1.Range("F20").Value = 479
$ws1.Range("F21").Value = 484
$ws1.Range("F22").Value = 343
$ws1.Range("F23").Value = 120
$ws1.Range("F25").Value = 1256
$ws1.Range("F26").Value = 105
$ws1.Range("F27").Value = 1457
$ws1.Range("F33").Value = 23
$ws1.Range("F34").Value = 70
$ws1.Range("F36").Value = 74
$ws1.Range("F37").Value = 3155
$ws1.Range("F39").Value = 48
$ws1.Range("F40").Value = 166
$ws1.Range("F41").Value = 50
$ws1.Range("F42").Value = 1149

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 25

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 5662
$ws4.Range("F4").Value  = 7670
$ws4.Range("F10").Value = 40
$ws4.Range("F11").Value = 4446
$ws4.Range("F15").Value = 3001
$ws4.Range("F19").Value = 551
$ws4.Range("F20").Value = 479
$ws4.Range("F21").Value = 484
$ws4.Range("F23").Value = 343
$ws4.Range("F24").Value = 120
$ws4.Range("F26").Value = 1256
$ws4.Range("F27").Value = 105
$ws4.Range("F28").Value = 1457
$ws4.Range("F34").Value = 23
$ws4.Range("F35").Value = 70
$ws4.Range("F37").Value = 74
$ws4.Range("F38").Value = 3155
$ws4.Range("F39").Value = 25
$ws4.Range("F41").Value = 48
$ws4.Range("F42").Value = 166
$ws4.Range("F43").Value = 50
$ws4.Range("F44").Value = 1149

$wb.Save()
